$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cell D1 with the same text and formatting as the other headers
$ws.Range("D1").Value = "Tipo"
$ws.Range("A1").Copy()
$ws.Range("D1").PasteSpecial(-4122)

# Update the existing model-evaluation numbers
$ws.Range("B2").Value = 0.501852495355711
$ws.Range("C2").Value = 0.9900070510728393

# Add the new "Tipo" data value
$ws.Range("D2").Value = "single"
